$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

$ALC.Cells.Item(15, 8).Value = 1789368.5
$ALC.Cells.Item(15, 9).Value = 1789368.5
$ALC.Cells.Item(15, 11).Value = 5368105.5
$ALC.Cells.Item(15, 13).Value = -5367936.5
$ALC.Cells.Item(32, 8).Value = 2334.35
$ALC.Cells.Item(32, 10).Value = 2192.1428
$ALC.Cells.Item(32, 12).Value = 2192.1428
$ALC.Cells.Item(32, 14).Value = -2844.1428
$ALC.Cells.Item(40, 8).Value = 2372.842
$ALC.Cells.Item(40, 9).Value = 1620
$ALC.Cells.Item(40, 10).Value = 2641.7144
$ALC.Cells.Item(40, 11).Value = 1620
$ALC.Cells.Item(40, 12).Value = 2641.7144
$ALC.Cells.Item(40, 13).Value = -1445
$ALC.Cells.Item(40, 14).Value = -2991.7144
$ALC.Cells.Item(64, 8).Value = 3817.7334
$ALC.Cells.Item(64, 9).Value = 3562.5
$ALC.Cells.Item(64, 10).Value = 3910.5454
$ALC.Cells.Item(64, 11).Value = 3562.5
$ALC.Cells.Item(64, 12).Value = 3910.5454
$ALC.Cells.Item(64, 13).Value = -3314.5
$ALC.Cells.Item(64, 14).Value = -4406.5454
$ALC.Cells.Item(67, 8).Value = 3817.7334
$ALC.Cells.Item(67, 9).Value = 3562.5
$ALC.Cells.Item(67, 10).Value = 3910.5454
$ALC.Cells.Item(67, 11).Value = 3562.5
$ALC.Cells.Item(67, 12).Value = 3910.5454
$ALC.Cells.Item(67, 13).Value = -2704.5
$ALC.Cells.Item(67, 14).Value = -5626.5454
$ALC.Cells.Item(74, 8).Value = 4603.75
$ALC.Cells.Item(74, 9).Value = 4346.6665
$ALC.Cells.Item(74, 10).Value = 5375
$ALC.Cells.Item(74, 11).Value = 4346.6665
$ALC.Cells.Item(74, 12).Value = 5375
$ALC.Cells.Item(74, 13).Value = -3410.6665
$ALC.Cells.Item(74, 14).Value = -7247
$ALC.Cells.Item(77, 8).Value = 4603.75
$ALC.Cells.Item(77, 9).Value = 4346.6665
$ALC.Cells.Item(77, 10).Value = 5375
$ALC.Cells.Item(77, 11).Value = 21733.3325
$ALC.Cells.Item(77, 12).Value = 26875
$ALC.Cells.Item(77, 13).Value = -17053.3325
$ALC.Cells.Item(77, 14).Value = -36235
$ALC.Cells.Item(100, 8).Value = 2192
$ALC.Cells.Item(100, 9).Value = 1675.7142
$ALC.Cells.Item(100, 10).Value = 3396.6667
$ALC.Cells.Item(100, 11).Value = 1675.7142
$ALC.Cells.Item(100, 12).Value = 3396.6667
$ALC.Cells.Item(100, 13).Value = -1134.7142
$ALC.Cells.Item(100, 14).Value = -4478.6667
$ALC.Cells.Item(107, 8).Value = 622
$ALC.Cells.Item(107, 9).Value = 485.33334
$ALC.Cells.Item(107, 10).Value = 950
$ALC.Cells.Item(107, 11).Value = 485.33334
$ALC.Cells.Item(107, 12).Value = 950
$ALC.Cells.Item(107, 13).Value = 1434.66666
$ALC.Cells.Item(107, 14).Value = -4790
$ALC.Cells.Item(113, 8).Value = 6043.091
$ALC.Cells.Item(113, 9).Value = 6722.7144
$ALC.Cells.Item(113, 10).Value = 4853.75
$ALC.Cells.Item(113, 11).Value = 6722.7144
$ALC.Cells.Item(113, 12).Value = 4853.75
$ALC.Cells.Item(113, 13).Value = -3468.7144
$ALC.Cells.Item(113, 14).Value = -11361.75
$ALC.Cells.Item(116, 8).Value = 689710.3
$ALC.Cells.Item(116, 9).Value = 2683.2
$ALC.Cells.Item(116, 10).Value = 1548494.1
$ALC.Cells.Item(116, 11).Value = 2683.2
$ALC.Cells.Item(116, 12).Value = 1548494.1
$ALC.Cells.Item(116, 13).Value = 758.8000000000002
$ALC.Cells.Item(116, 14).Value = -1555378.1
$ARM.Cells.Item(2, 8).Value = 3868.3125
$ARM.Cells.Item(2, 9).Value = 3076.3845
$ARM.Cells.Item(2, 11).Value = 3076.3845
$ARM.Cells.Item(2, 13).Value = -2963.3845
$ARM.Cells.Item(39, 8).Value = 0
$ARM.Cells.Item(39, 9).Value = 0
$ARM.Cells.Item(39, 11).Value = 0
$ARM.Cells.Item(39, 13).ClearContents()
$ARM.Cells.Item(45, 8).Value = 3344.1428
$ARM.Cells.Item(45, 9).Value = 1219.9
$ARM.Cells.Item(45, 10).Value = 5275.273
$ARM.Cells.Item(45, 11).Value = 1219.9
$ARM.Cells.Item(45, 12).Value = 5275.273
$ARM.Cells.Item(45, 13).Value = -842.9000000000001
$ARM.Cells.Item(45, 14).Value = -6029.273
$ARM.Cells.Item(57, 8).Value = 19500
$ARM.Cells.Item(57, 9).Value = 19500
$ARM.Cells.Item(57, 11).Value = 19500
$ARM.Cells.Item(57, 13).Value = -19016
$ARM.Cells.Item(116, 8).Value = 3868.3125
$ARM.Cells.Item(116, 9).Value = 3076.3845
$ARM.Cells.Item(116, 11).Value = 3076.3845
$ARM.Cells.Item(116, 13).Value = -782.3845000000001
$BSM.Cells.Item(3, 8).Value = 3868.3125
$BSM.Cells.Item(3, 9).Value = 3076.3845
$BSM.Cells.Item(3, 11).Value = 3076.3845
$BSM.Cells.Item(3, 13).Value = -2962.3845
$BSM.Cells.Item(113, 8).Value = 4885
$BSM.Cells.Item(113, 9).Value = 4885
$BSM.Cells.Item(113, 11).Value = 4885
$BSM.Cells.Item(113, 13).Value = -2715
$CRP.Cells.Item(7, 8).Value = 13078.75
$CRP.Cells.Item(7, 9).Value = 100000
$CRP.Cells.Item(7, 10).Value = 661.4286
$CRP.Cells.Item(7, 11).Value = 100000
$CRP.Cells.Item(7, 12).Value = 661.4286
$CRP.Cells.Item(7, 13).Value = -99887
$CRP.Cells.Item(7, 14).Value = -887.4286
$CRP.Cells.Item(22, 8).Value = 495.8
$CRP.Cells.Item(22, 9).Value = 262.6
$CRP.Cells.Item(22, 10).Value = 729
$CRP.Cells.Item(22, 11).Value = 262.6
$CRP.Cells.Item(22, 12).Value = 729
$CRP.Cells.Item(22, 13).Value = 87.39999999999998
$CRP.Cells.Item(22, 14).Value = -1429
$CRP.Cells.Item(36, 8).Value = 1500
$CRP.Cells.Item(36, 9).Value = 1500
$CRP.Cells.Item(36, 10).Value = 0
$CRP.Cells.Item(36, 11).Value = 1500
$CRP.Cells.Item(36, 12).Value = 0
$CRP.Cells.Item(36, 13).Value = -1112
$CRP.Cells.Item(36, 14).ClearContents()
$CRP.Cells.Item(38, 8).Value = 4397.5
$CRP.Cells.Item(38, 9).Value = 3196.6667
$CRP.Cells.Item(38, 10).Value = 8000
$CRP.Cells.Item(38, 11).Value = 3196.6667
$CRP.Cells.Item(38, 12).Value = 8000
$CRP.Cells.Item(38, 13).Value = -2819.6667
$CRP.Cells.Item(38, 14).Value = -8754
$CRP.Cells.Item(40, 8).Value = 1500
$CRP.Cells.Item(40, 9).Value = 1500
$CRP.Cells.Item(40, 10).Value = 0
$CRP.Cells.Item(40, 11).Value = 1500
$CRP.Cells.Item(40, 12).Value = 0
$CRP.Cells.Item(40, 13).Value = -1340
$CRP.Cells.Item(40, 14).ClearContents()
$CRP.Cells.Item(42, 8).Value = 0
$CRP.Cells.Item(42, 10).Value = 0
$CRP.Cells.Item(42, 12).Value = 0
$CRP.Cells.Item(42, 14).ClearContents()
$CRP.Cells.Item(46, 8).Value = 4397.5
$CRP.Cells.Item(46, 9).Value = 3196.6667
$CRP.Cells.Item(46, 10).Value = 8000
$CRP.Cells.Item(46, 11).Value = 3196.6667
$CRP.Cells.Item(46, 12).Value = 8000
$CRP.Cells.Item(46, 13).Value = -2985.6667
$CRP.Cells.Item(46, 14).Value = -8422
$CUL.Cells.Item(64, 8).Value = 2877.2727
$CUL.Cells.Item(64, 9).Value = 1650
$CUL.Cells.Item(64, 10).Value = 3150
$CUL.Cells.Item(64, 11).Value = 4950
$CUL.Cells.Item(64, 12).Value = 9450
$CUL.Cells.Item(64, 13).Value = -4680
$CUL.Cells.Item(64, 14).Value = -9990
$CUL.Cells.Item(67, 8).Value = 2877.2727
$CUL.Cells.Item(67, 9).Value = 1650
$CUL.Cells.Item(67, 10).Value = 3150
$CUL.Cells.Item(67, 11).Value = 4950
$CUL.Cells.Item(67, 12).Value = 9450
$CUL.Cells.Item(67, 13).Value = -4014
$CUL.Cells.Item(67, 14).Value = -11322
$GSM.Cells.Item(11, 8).Value = 7080555.5
$GSM.Cells.Item(11, 9).Value = 5955384.5
$GSM.Cells.Item(11, 10).Value = 10006000
$GSM.Cells.Item(11, 11).Value = 5955384.5
$GSM.Cells.Item(11, 12).Value = 10006000
$GSM.Cells.Item(11, 13).Value = -5955245.5
$GSM.Cells.Item(11, 14).Value = -10006278
$GSM.Cells.Item(46, 8).Value = 29996
$GSM.Cells.Item(46, 10).Value = 29996
$GSM.Cells.Item(46, 12).Value = 29996
$GSM.Cells.Item(46, 14).Value = -30308
$LTW.Cells.Item(107, 8).Value = 4385
$LTW.Cells.Item(107, 9).Value = 4385
$LTW.Cells.Item(107, 11).Value = 4385
$LTW.Cells.Item(107, 13).Value = -2465
$WVR.Cells.Item(126, 8).Value = 3006.6538
$WVR.Cells.Item(126, 9).Value = 3074.9048
$WVR.Cells.Item(126, 10).Value = 2720
$WVR.Cells.Item(126, 11).Value = 9224.714399999999
$WVR.Cells.Item(126, 12).Value = 8160
$WVR.Cells.Item(126, 13).Value = -6754.714399999999
$WVR.Cells.Item(126, 14).Value = -13100
